# Auto-generated edit script to apply cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.844.25"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "3.529.91"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'605.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'197.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.80%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.13%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("D11").Value = "'53.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "4.088.52"
$ws.Range("D15").Value = "'597.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "69.991.49"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "'19.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "'12.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "3.512.81"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'18.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.71%  "
$ws.Range("D23").Value = "'5.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.35%  "
$ws.Range("D24").Value = "'101.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.42%  "
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'9.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.16%  "
$ws.Range("D29").Value = "'33.44"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("B30").Value = "dogwifhat"
$ws.Range("C30").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D30").Value = "'4.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.04%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'7.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'12.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").Value = "'63.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "0.0₃0856"
$ws.Range("E35").Value = "  +10.52%  "
$ws.Range("D36").Value = "3.733.25"
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("D39").Value = "'3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "'36.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").Value = "'489.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.42%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'0.0455"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("D49").Value = "'8.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "'0.000252"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").Value = "'130.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
